$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: remove the "Meta description" paragraph that sits right after the
# title (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -gt 0) {
    $d.Paragraphs($metaIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# Step 2: locate the final "Prompt: DALLE..." paragraph and, right before it,
# insert a new bold paragraph containing the page title.
# ---------------------------------------------------------------------------
$dalleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Prompt: DALLE*") {
        $dalleIndex = $i
        break
    }
}

$dalleP = $d.Paragraphs($dalleIndex)
$insertPoint = $d.Range($dalleP.Range.Start, $dalleP.Range.Start)

# Build the new paragraph with the same run shape as the rest of the
# document (a leading empty run followed by a bold text run), using an
# OOXML fragment so the formatting is exact and not inherited from either
# neighbouring paragraph.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Eggomatic Free: A Unique Steampunk Slot Game</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml) | Out-Null

# InsertXML splits the DALLE paragraph's leading edge in two, leaving a
# stray empty paragraph between our new title paragraph and the DALLE
# paragraph; remove that stray paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "`r") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Step 3: replace the DALLE prompt text with the meta-description text
# (the run keeps its existing italic formatting).
# ---------------------------------------------------------------------------
$dalleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Prompt: DALLE*") {
        $dalleIndex = $i
        break
    }
}
$dalleP = $d.Paragraphs($dalleIndex)
$r = $d.Range($dalleP.Range.Start, $dalleP.Range.End)
$r.Text = "Read our review of Eggomatic, a top-notch online slot game with a unique steampunk theme and exciting bonus features. Play for free now."
